$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of quotes for 2025-09-25 (serial date 45925)
$row = 21

$ws.Cells.Item($row, 1).Value = 45925
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = "14,9067"
$ws.Cells.Item($row, 3).Value = "15,0772"
$ws.Cells.Item($row, 4).Value = "14,9067"
$ws.Cells.Item($row, 5).Value = "14,9067"
